# Apply updated benchmark numbers (source data changed after a bool-instead-of-int
# fix in the upstream sdql-rs code), then leave the selection where the author
# last left it before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("benchmarks")

# Row 5 (query #4): DuckDB time dropped 39 -> 37, sdql-rs time dropped 8 -> 6.
# Dependent "speedup"/"parallel vs seq" formulas recalc automatically.
$ws.Range("C5").Value = 37
$ws.Range("E5").Value = 6

# Row 22 (query #21): DuckDB time rose 671 -> 674.
$ws.Range("C22").Value = 674

# Row 23 (query #22): sdql-rs (parallel) time rose 2 -> 3.
$ws.Range("E23").Value = 3

# Leave the saved selection on A24, matching the author's last cursor position.
$ws.Range("A24").Select()
